# The workbook had an extra "contraseña nueva" column (K) that was removed
# by the author, shifting the "Ingresos" column (previously L) left into K.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "contraseña nueva" column (K). This removes its header
# cell and all of its data cells, shifts every column to the right of it
# (just "Ingresos" in L) one position to the left, and drops the now-unused
# shared string from the workbook when it is saved.
$ws.Columns("K").Delete()

# Match the author's final selection in the saved file.
[void]$ws.Range("C3").Select()
